{"js": "// 1. Title: \"Use Case \"Select a dock marked on map\"\" -> \"Use Case \"Select a dock on list\"\"\nconst titleResults = context.document.body.search(\"marked on map\", { matchCase: true });\ntitleResults.load(\"text\");\nawait context.sync();\ntitleResults.items[0].insertText(\"on list\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2. \"B\u01b0\u1edbc 1: Kh\u00e1ch h\u00e0ng nh\u1ea5n v\u00e0o c\u00e1c \u0111i\u1ec3m \u0111\u00e3 \u0111\u01b0\u1ee3c \u0111\u00e1nh d\u1ea5u tr\u00ean b\u1ea3n \u0111\u1ed3.\"\n//    -> \"B\u01b0\u1edbc 1: Kh\u00e1ch h\u00e0ng l\u1ef1a ch\u1ecdn b\u00e3i xe tr\u00ean danh s\u00e1ch.\"\nconst step1Results = context.document.body.search(\n  \"nh\u1ea5n v\u00e0o c\u00e1c \u0111i\u1ec3m \u0111\u00e3 \u0111\u01b0\u1ee3c \u0111\u00e1nh d\u1ea5u tr\u00ean b\u1ea3n \u0111\u1ed3\",\n  { matchCase: true }\n);\nstep1Results.load(\"text\");\nawait context.sync();\nstep1Results.items[0].insertText(\"l\u1ef1a ch\u1ecdn b\u00e3i xe tr\u00ean danh s\u00e1ch\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3. \"B\u01b0\u1edbc 1: Kh\u00e1ch h\u00e0ng nh\u1ea5n v\u00e0o thanh t\u00ecm ki\u1ebfm ph\u00eda tr\u00ean b\u1ea3n \u0111\u1ed3 \u0111\u1ec3 nh\u1eadp. \"\n//    -> \"B\u01b0\u1edbc 1: Kh\u00e1ch h\u00e0ng nh\u1ea5n v\u00e0o thanh t\u00ecm ki\u1ebfm ph\u00eda tr\u00ean danh s\u00e1ch \u0111\u1ec3 nh\u1eadp. \"\nconst searchBarResults = context.document.body.search(\n  \"t\u00ecm ki\u1ebfm ph\u00eda tr\u00ean b\u1ea3n \u0111\u1ed3 \u0111\u1ec3 nh\u1eadp\",\n  { matchCase: true }\n);\nsearchBarResults.load(\"text\");\nawait context.sync();\nsearchBarResults.items[0].insertText(\"t\u00ecm ki\u1ebfm ph\u00eda tr\u00ean danh s\u00e1ch \u0111\u1ec3 nh\u1eadp\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 4. Merge the split runs of \"B\u01b0\u1edbc 4: H\u1ec7 th\u1ed1ng s\u1ebd hi\u1ec7n l\u00ean c\u00e1c th\u00f4ng tin v\u1ec1 b\u00e3i xe m\u00e0 kh\u00e1ch h\u00e0ng \u0111\u00e3 ch\u1ecdn.\"\n//    back into a single run (no visible text change).\nconst step4Text = \"B\u01b0\u1edbc 4: H\u1ec7 th\u1ed1ng s\u1ebd hi\u1ec7n l\u00ean c\u00e1c th\u00f4ng tin v\u1ec1 b\u00e3i xe m\u00e0 kh\u00e1ch h\u00e0ng \u0111\u00e3 ch\u1ecdn.\";\nconst step4Results = context.document.body.search(step4Text, { matchCase: true });\nstep4Results.load(\"text\");\nawait context.sync();\nstep4Results.items[0].insertText(step4Text, Word.InsertLocation.replace);\nawait context.sync();\n\n// 5. Move the \"_GoBack\" bookmark from the empty paragraph (after step 4) to the\n//    end of the document's last content run (\"...xem th\u00f4ng tin.\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst bodyParagraphs = context.document.body.paragraphs;\nbodyParagraphs.load(\"items\");\nawait context.sync();\nconst lastParagraph = bodyParagraphs.items[bodyParagraphs.items.length - 2];\nconst endRange = lastParagraph.getRange(Word.RangeLocation.end);\nendRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($searchText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n\n# 1. Title: \"Use Case \"Select a dock marked on map\"\" -> \"Use Case \"Select a dock on list\"\"\nReplace-Text \"marked on map\" \"on list\"\n\n# 2. \"B\u01b0\u1edbc 1: Kh\u00e1ch h\u00e0ng nh\u1ea5n v\u00e0o c\u00e1c \u0111i\u1ec3m \u0111\u00e3 \u0111\u01b0\u1ee3c \u0111\u00e1nh d\u1ea5u tr\u00ean b\u1ea3n \u0111\u1ed3.\"\n#    -> \"B\u01b0\u1edbc 1: Kh\u00e1ch h\u00e0ng l\u1ef1a ch\u1ecdn b\u00e3i xe tr\u00ean danh s\u00e1ch.\"\nReplace-Text \"nh\u1ea5n v\u00e0o c\u00e1c \u0111i\u1ec3m \u0111\u00e3 \u0111\u01b0\u1ee3c \u0111\u00e1nh d\u1ea5u tr\u00ean b\u1ea3n \u0111\u1ed3\" \"l\u1ef1a ch\u1ecdn b\u00e3i xe tr\u00ean danh s\u00e1ch\"\n\n# 3. \"B\u01b0\u1edbc 1: Kh\u00e1ch h\u00e0ng nh\u1ea5n v\u00e0o thanh t\u00ecm ki\u1ebfm ph\u00eda tr\u00ean b\u1ea3n \u0111\u1ed3 \u0111\u1ec3 nh\u1eadp. \"\n#    -> \"B\u01b0\u1edbc 1: Kh\u00e1ch h\u00e0ng nh\u1ea5n v\u00e0o thanh t\u00ecm ki\u1ebfm ph\u00eda tr\u00ean danh s\u00e1ch \u0111\u1ec3 nh\u1eadp. \"\nReplace-Text \"t\u00ecm ki\u1ebfm ph\u00eda tr\u00ean b\u1ea3n \u0111\u1ed3 \u0111\u1ec3 nh\u1eadp\" \"t\u00ecm ki\u1ebfm ph\u00eda tr\u00ean danh s\u00e1ch \u0111\u1ec3 nh\u1eadp\"\n\n# 4. Merge the split runs of \"B\u01b0\u1edbc 4: H\u1ec7 th\u1ed1ng s\u1ebd hi\u1ec7n l\u00ean c\u00e1c th\u00f4ng tin v\u1ec1 b\u00e3i xe m\u00e0 kh\u00e1ch h\u00e0ng \u0111\u00e3 ch\u1ecdn.\"\n#    back into a single run (no visible text change).\n$step4Text = \"B\u01b0\u1edbc 4: H\u1ec7 th\u1ed1ng s\u1ebd hi\u1ec7n l\u00ean c\u00e1c th\u00f4ng tin v\u1ec1 b\u00e3i xe m\u00e0 kh\u00e1ch h\u00e0ng \u0111\u00e3 ch\u1ecdn.\"\nReplace-Text $step4Text $step4Text\n\n# 5. Move the \"_GoBack\" bookmark from the empty paragraph (after step 4) to the\n#    end of the document's last content run (\"...xem th\u00f4ng tin.\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"xem th\u00f4ng tin.\"\n$find2.Execute() | Out-Null\n$afterPhraseStart = $d.Content.End  # placeholder, recomputed below\n\n# Locate the end offset of \"xem th\u00f4ng tin.\" via Find, then work around an\n# engine quirk where a zero-length range exactly at a paragraph's text-end\n# offset gets mis-resolved by Bookmarks.Add: insert a 1-char placeholder\n# right after the target spot (so the bookmark position is no longer a\n# paragraph-end offset), add the bookmark there, then delete the\n# placeholder again. The bookmark stays correctly anchored in place.\n$markerRange = $d.Content\n$markerFind = $markerRange.Find\n$markerFind.ClearFormatting()\n$markerFind.Text = \"xem th\u00f4ng tin.\"\n$markerFind.Execute() | Out-Null\n$targetPos = $markerRange.End\n\n$placeholderRange = $d.Range($targetPos, $targetPos)\n$placeholderRange.InsertAfter(\"X\")\n\n$bmRange = $d.Range($targetPos, $targetPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n$d.Range($targetPos, $targetPos + 1).Delete()\n"}
